$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308, shifting existing rows 308:379 down to 309:380
$ws.Rows("308:308").Insert()

# Fill in the new row 308 with the new record's data
$ws.Range("A308").Value = 4
$ws.Range("B308").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C308").Value = "Los Lagos"
$ws.Range("D308").Value = 45015
$ws.Range("E308").Value = 10
$ws.Range("F308").Value = 100112024
$ws.Range("G308").Value = "Choclo"
$ws.Range("H308").Value = "Dulce o Americano"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 10000
$ws.Range("K308").Value = 300
$ws.Range("L308").Value = 300
$ws.Range("M308").Value = 300
$ws.Range("N308").Value = "$/unidad"
$ws.Range("O308").Value = "Región de La Araucanía"
$ws.Range("P308").Value = 300
$ws.Range("Q308").Value = 1
$ws.Range("R308").Value = "Hortaliza"
